$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row at row 33 (everything from row 33 downward shifts down by one),
# matching the new "READCOUNT / 读取次数" entry added to the T_ARTICLE field table.
$ws.Rows.Item(33).Insert(-4121)

# Copy the formatting of the row above (the last existing T_ARTICLE data row, now
# row 32) into the freshly inserted row 33 so the cell styles line up exactly with
# the rest of the table (border, font, alignment, wrap, etc.).
$ws.Range("A32:F32").Copy()
$ws.Range("A33:F33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row with the READCOUNT field describing the read-count column.
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "READCOUNT"
$ws.Range("C33").Value = "INT"
$ws.Range("D33").Value = "读取次数"
$ws.Range("E33").Value = $null
$ws.Range("F33").Value = $null

# Update the view state to match: scrolled back up a bit, with D33 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D33").Select()
